$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 15 updates (group: dsa_discussion_group)
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = "2026-02-19T07:33:19.034492+00:00"
$ws.Range("H15").Value = 8
$ws.Range("L15").Value = "[136, 155, 150, 158, 142, 140, 151, 137]"
